# Auto-generated edit script applying the cryptos list update (GitHub Actions cron refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.999.85'
$ws.Range('E2').Value = '''  -0.33%  '
$ws.Range('D3').Value = '''3.319.57'
$ws.Range('E3').Value = '''  +0.32%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '''  +0.07%  '
$ws.Range('D5').Value = '''581.05'
$ws.Range('E5').Value = '''  +3.94%  '
$ws.Range('D6').Value = '''184.71'
$ws.Range('E6').Value = '''  -2.71%  '
$ws.Range('E7').Value = '''  +0.09%  '
$ws.Range('D8').Value = '''3.313.44'
$ws.Range('E8').Value = '''  +0.35%  '
$ws.Range('D9').Value = '''0.573'
$ws.Range('E9').Value = '''  -2.76%  '
$ws.Range('D10').Value = '''0.179'
$ws.Range('E10').Value = '''  -3.01%  '
$ws.Range('D11').Value = '''0.575'
$ws.Range('E11').Value = '''  -2.27%  '
$ws.Range('D12').Value = '''46.86'
$ws.Range('E12').Value = '''  -2.41%  '
$ws.Range('D13').Value = '''0.0000266'
$ws.Range('E13').Value = '''  -1.76%  '
$ws.Range('D14').Value = '''666.26'
$ws.Range('E14').Value = '''  +8.92%  '
$ws.Range('D15').Value = '''3.855.44'
$ws.Range('E15').Value = '''  +0.38%  '
$ws.Range('D16').Value = '''8.42'
$ws.Range('E16').Value = '''  -3.36%  '
$ws.Range('D17').Value = '''66.164.44'
$ws.Range('E17').Value = '''  -0.13%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '''0.117'
$ws.Range('E18').Value = '''  -0.29%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''17.84'
$ws.Range('E19').Value = '''  -1.36%  '
$ws.Range('D20').Value = '''3.325.36'
$ws.Range('E20').Value = '''  +0.48%  '
$ws.Range('D21').Value = '''11.05'
$ws.Range('E21').Value = '''  -0.62%  '
$ws.Range('D22').Value = '''0.893'
$ws.Range('E22').Value = '''  -2.20%  '
$ws.Range('D23').Value = '''17.81'
$ws.Range('E23').Value = '''  -3.51%  '
$ws.Range('E24').Value = '''  -0.60%  '
$ws.Range('D25').Value = '''5.01'
$ws.Range('E25').Value = '''  -1.81%  '
$ws.Range('E26').Value = '''  -0.56%  '
$ws.Range('D28').Value = '''9.48'
$ws.Range('E28').Value = '''  -3.51%  '
$ws.Range('D29').Value = '''31.26'
$ws.Range('E29').Value = '''  +3.02%  '
$ws.Range('D30').Value = '''8.43'
$ws.Range('E30').Value = '''  -2.77%  '
$ws.Range('D31').Value = '''6.65'
$ws.Range('E31').Value = '''  -1.90%  '
$ws.Range('D32').Value = '''590.44'
$ws.Range('E32').Value = '''  +3.08%  '
$ws.Range('D33').Value = '''3.83'
$ws.Range('D34').Value = '''10.96'
$ws.Range('E34').Value = '''  -1.45%  '
$ws.Range('D35').Value = '''0.105'
$ws.Range('E35').Value = '''  -0.40%  '
$ws.Range('D36').Value = '''3.833.17'
$ws.Range('E37').Value = '''  +0.17%  '
$ws.Range('D38').Value = '''55.77'
$ws.Range('E38').Value = '''  -2.52%  '
$ws.Range('D39').Value = '''2.66'
$ws.Range('E39').Value = '''  -2.10%  '
$ws.Range('D40').Value = '''32.75'
$ws.Range('E40').Value = '''  -4.35%  '
$ws.Range('D41').Value = '''0.0₃0695'
$ws.Range('E41').Value = '''  -5.07%  '
$ws.Range('E42').Value = '''  -3.64%  '
$ws.Range('D43').Value = '''3.17'
$ws.Range('E43').Value = '''  -4.80%  '
$ws.Range('D44').Value = '''3.40'
$ws.Range('E44').Value = '''  +4.17%  '
$ws.Range('D45').Value = '''0.334'
$ws.Range('E45').Value = '''  -2.10%  '
$ws.Range('D46').Value = '''0.0410'
$ws.Range('E46').Value = '''  -3.89%  '
$ws.Range('D47').Value = '''3.02'
$ws.Range('E47').Value = '''  -13.59%  '
$ws.Range('D48').Value = '''0.127'
$ws.Range('E48').Value = '''  -1.78%  '
$ws.Range('E49').Value = '''  +0.32%  '
$ws.Range('D50').Value = '''2.54'
$ws.Range('E50').Value = '''  -2.45%  '
$ws.Range('D51').Value = '''130.17'
$ws.Range('E51').Value = '''  +5.81%  '
